$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllTablesWithGaps")

$oldText1 = '"Part Number - Can be found on the top right position of the page"'
$newText1 = '"Part Number - Can be found on the center right position of the page"'

$oldText2 = '"Duxford Range Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the bottom right position of the page"'
$newText2 = '"Multi-Wattage Tri-Colour and Single Colour 4000K Retrofit Gear Trays - Can be found on the middle right position of the page"'

foreach ($r in 2..13) {
    $cell = $ws.Range("G$r")
    if ($cell.Text -eq $oldText1) {
        $cell.Value = $newText1
    }
}

foreach ($r in 17..28) {
    $cell = $ws.Range("G$r")
    if ($cell.Text -eq $oldText2) {
        $cell.Value = $newText2
    }
}
